$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cells touched below hold plain numeric-looking / percentage text
# values (coin prices, % changes) that must remain TEXT, exactly as
# originally authored, instead of being reinterpreted by Excel as
# numbers (which would drop trailing zeros / switch to scientific
# notation, e.g. "0.0000159" -> 1.59E-05). So for every cell we are
# about to write, force Text format first, then set its value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.862.90"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.768.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.14%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.74"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +11.28%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.58%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.996"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.788.02"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.15%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.97%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.398"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.20%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.261.88"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.53"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.809.76"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000159"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +7.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.785.75"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.33"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.62%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "368.06"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.08"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.556"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.92%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.30"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.71%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.71"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0963"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +15.83%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.04"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.35"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.27"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +11.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "172.93"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.82"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.11%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.07"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +8.35%  "

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +7.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.84"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.04"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "344.14"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.40%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.24"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +14.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.98"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.45"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.94"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +8.01%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.653"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.69%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "138.66"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.01%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.183.24"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.28%  "

